# Applies the "build site" content update for LOQ4222.xlsx
# - Updates the activation date
# - Rewrites Objetivos / adds the missing Objectives (English) text
# - Updates the responsible professor
# - Rewrites Programa resumido / adds missing Short syllabus (English) text
# - Rewrites Programa / adds missing Syllabus (English) text
# - Rewrites Método, Critério, Norma de recuperação and Bibliografia

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

function Set-TextValue($range, [string]$value) {
    # Force the cell to remain a text value even when the content looks
    # like a date (e.g. "01/01/2021"), which Excel would otherwise
    # auto-convert into a date serial number.
    $range.NumberFormat = "@"
    $range.Value = $value
}

function Copy-CellFormat($sourceRange, $targetRange) {
    $sourceRange.Copy()
    $targetRange.PasteSpecial($xlPasteFormats)
}

# --- Ativação: 01/01/2018 -> 01/01/2021 -----------------------------------
Set-TextValue $ws.Range("B8") "01/01/2021"
Copy-CellFormat $ws.Range("B7") $ws.Range("B8")
Set-TextValue $ws.Range("C8") "01/01/2021"
Copy-CellFormat $ws.Range("C7") $ws.Range("C8")

# --- Objetivos: (row 10) new Portuguese text -------------------------------
$objetivosPt = "Apresentar noções de Matemática Financeira, Gestão Financeiras e Engenharia Econômica aos alunos, capacitando-os para construir e analisar fluxos de caixa de projetos e empreendimentos, discutir os principais aspectos da gestão financeira das empresas industriais, comerciais e de serviços e analisar e propor estratégias de gestão financeira relacionadas às estratégias de mercado e de produção."
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# --- Objectives: (row 11) previously empty B/C, now filled (English) ------
$objetivosEn = "To present notions of Financial Mathematics, Financial Management and Economic Engineering to students, enabling them to build and analyze cash flows from projects and enterprises, discuss the main aspects of financial management in industrial, commercial and service companies and analyze and propose strategies for financial management related to market and production strategies."
Copy-CellFormat $ws.Range("B10") $ws.Range("B11")
$ws.Range("B11").Value = $objetivosEn
Copy-CellFormat $ws.Range("C10") $ws.Range("C11")
$ws.Range("C11").Value = $objetivosEn

# --- Docentes responsáveis: (row 13) new professor -------------------------
$docente = "11079086 - Herlandí de Souza Andrade"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# --- Programa resumido: (row 14) minor text change --------------------------
$resumido = "1. Engenharia Econômica. 2 – Finanças"
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# --- Short syllabus: (row 15) previously empty B/C, now filled (English) --
$shortSyllabusEn = "1. Economic Engineering. 2 – Finance"
Copy-CellFormat $ws.Range("B14") $ws.Range("B15")
$ws.Range("B15").Value = $shortSyllabusEn
Copy-CellFormat $ws.Range("C14") $ws.Range("C15")
$ws.Range("C15").Value = $shortSyllabusEn

# --- Programa: (row 16) updated Portuguese syllabus text --------------------
$programa = "1. Engenharia Econômica: Variável tempo: juros simples, juros compostos; Métodos de amortização; Equivalência de métodos; Métodos de Decisão; Renovação e substituição de equipamentos; Depreciação; Análise de Projetos, Riscos em projetos; Estimativa do custo de capital próprio (CAPM) e WACC.2. Finanças: O ciclo da produção e o ciclo do capital; Análise de Índices; Fontes de Financiamento, Alavancagem; Capital de Giro; Custo de Capital; Ações, Política de Dividendos; Financiamento de Longo Prazo, Corporate Finance/Project Finance; EVA e MVA."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Syllabus: (row 17) previously empty B/C, now filled (English) --------
$syllabusEn = "1. Economic Engineering: Time variable: simple interest, compound interest; Amortization methods; Equivalence of methods; Decision Methods; Renovation and replacement of equipment; Depreciation; Project Analysis, Project Risks; Estimated cost of equity (CAPM) and WACC.2. Finance: The production cycle and the capital cycle; Index Analysis; Financing Sources, Leverage; Working capital; Capital cost; Shares, Dividend Policy; Long Term Financing, Corporate Finance/Project Finance; EVA and MVA."
Copy-CellFormat $ws.Range("B16") $ws.Range("B17")
$ws.Range("B17").Value = $syllabusEn
Copy-CellFormat $ws.Range("C16") $ws.Range("C17")
$ws.Range("C17").Value = $syllabusEn

# --- Método: (row 19) updated text ------------------------------------------
$metodo = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Critério: (row 20) updated text ----------------------------------------
$criterio = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Norma de recuperação: (row 21) minor text change -----------------------
$norma = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- Bibliografia: (row 22) entirely new bibliography text ------------------
$biblio = "GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.ASSAF NETO, A. E LIMA, F. G. 3 ed. CURSO DE ADMINISTRAÇÃO FINANCEIRA. São Paulo: Atlas, 2014MARIANO, F.; MENESES, A. Curso De Administração Financeira. São Paulo: Método, 2012.MORANTE, A. S. Análise das Demonstrações Financeiras. 2 ed. São Paulo: Atlas, 2009.NEWNAN, D. G.; LAVELLE, J. P. Fundamentos de Engenharia Econômica. São Paulo: LTC,2000.KOPITTKE, B. H.; CASAROTTO FILHO, N. ANÁLISE DE INVESTIMENTOS: Matemática Financeira, Engenharia Econômica, Estratégia Empresarial. 11 ed. São Paulo: Atlas, 2010.HOJI, M.; LUZ, A. E. Gestão Financeira Econômica: Didática, Objetiva e Prática. São Paulo: Atlas, 2019.GOMES, J. M. Elaboração e Análise De Viabilidade Econômica De Projetos. São Paulo: Atlas, 2013.OLIVO, R. L. F. Análise de Investimentos. Campinas: Alínea, 2011.ALMEIDA, J. T. S. Matemática Financeira. Rio de Janeiro: LTC, 2016.EHRLICH, Pierre Jacques. Engenharia Econômica. São Paulo: Editora Atlas, 2005.HIRSCHFELD, Henrique. Engenharia econômica e análise de custos. 7. ed. São Paulo: Atlas, 2007.MOTTA, Regis da Rocha; CALÔBA, Guilherme Marques. Análise de Investimentos. São Paulo: Atlas 2002.SANVICENTE, A. Z. Administração Financeira. São Paulo: Editora Atlas, 2007.VAN HORNE, J. C. Política e Administração Financeira. Rio de Janeiro: Livros Técnicos e Científicos, 1974.WESTON, J. F.; BRIGHAM, E. F. Administração Financeira de Empresas. São Paulo: Editora Interamericana, 2000."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio

$ws.Application.CutCopyMode = $false
